$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.7376376588883126
$ws.Range("J2").Value = 0.7376376588883125
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.022792
$ws.Range("N2").Value = 0.06837600000000001
$ws.Range("O2").Value = 0.001916327914826657
$ws.Range("P2").Value = 0.001916327914826657
$ws.Range("Q2").Value = 0.01327804180266667
$ws.Range("R2").Value = 0.119502376224
$ws.Range("S2").Value = 0.001413555636755057
$ws.Range("T2").Value = 0.001413555636755057
$ws.Range("I3").Value = 0.7376376588883126
$ws.Range("J3").Value = 0.7376376588883125
$ws.Range("O3").Value = 0.3701235913233977
$ws.Range("P3").Value = 0.3701235913233977
$ws.Range("S3").Value = 0.2730170994031256
$ws.Range("T3").Value = 0.2730170994031256
$ws.Range("I4").Value = 0.7376376588883126
$ws.Range("J4").Value = 0.7376376588883125
$ws.Range("M4").Value = 7.468693666666667
$ws.Range("N4").Value = 22.406081
$ws.Range("O4").Value = 0.6279600807617757
$ws.Range("P4").Value = 0.6279600807617757
$ws.Range("Q4").Value = 4.351071723293778
$ws.Range("R4").Value = 39.159645509644
$ws.Range("S4").Value = 0.4632070038484319
$ws.Range("T4").Value = 0.4632070038484318
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2072096666666667
$ws.Range("H5").Value = 0.621629
$ws.Range("I5").Value = 0.2623623411116874
$ws.Range("J5").Value = 0.2623623411116874
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.022792
$ws.Range("N5").Value = 0.06837600000000001
$ws.Range("O5").Value = 0.001916327914826657
$ws.Range("P5").Value = 0.001916327914826657
$ws.Range("Q5").Value = 0.004722722722666667
$ws.Range("R5").Value = 0.042504504504
$ws.Range("S5").Value = 0.0005027722780716001
$ws.Range("T5").Value = 0.0005027722780716001
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.2072096666666667
$ws.Range("H6").Value = 0.621629
$ws.Range("I6").Value = 0.2623623411116874
$ws.Range("J6").Value = 0.2623623411116874
$ws.Range("O6").Value = 0.3701235913233977
$ws.Range("P6").Value = 0.3701235913233977
$ws.Range("Q6").Value = 0.9121565685151111
$ws.Range("R6").Value = 8.209409116635999
$ws.Range("S6").Value = 0.09710649192027207
$ws.Range("T6").Value = 0.09710649192027207
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.2072096666666667
$ws.Range("H7").Value = 0.621629
$ws.Range("I7").Value = 0.2623623411116874
$ws.Range("J7").Value = 0.2623623411116874
$ws.Range("M7").Value = 7.468693666666667
$ws.Range("N7").Value = 22.406081
$ws.Range("O7").Value = 0.6279600807617757
$ws.Range("P7").Value = 0.6279600807617757
$ws.Range("Q7").Value = 1.547585525105444
$ws.Range("R7").Value = 13.928269725949
$ws.Range("S7").Value = 0.1647530769133438
$ws.Range("T7").Value = 0.1647530769133438
